$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily portfolio-update row (row 92) below the existing
# data (rows 1..91). Column A holds the date as a literal text label (the
# sheet stores dates as plain text, not Excel date serials) while columns
# B/C/D hold the numeric closing values.
$row = 92

# Use a leading apostrophe so Excel stores the date string as literal text
# instead of auto-converting it into a date serial number, matching how
# every other row in column A is stored.
$ws.Cells.Item($row, 1).Value = "'2025-11-15"
$ws.Cells.Item($row, 2).Value = 57.68000030517578
$ws.Cells.Item($row, 3).Value = 391.2000122070312
$ws.Cells.Item($row, 4).Value = 303.75
